$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need the
# NumberFormat forced to Text ("@") BEFORE assigning the value,
# otherwise Excel will auto-convert them into numeric values and
# lose formatting such as trailing zeros (e.g. "6.76" -> 6.76 float).
# Values containing two decimal separators (e.g. "69.025.11") are
# never auto-parsed by Excel, so they do not need this treatment.

$ws.Range("D2").Value = '69.025.11'
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").Value = '3.770.68'
$ws.Range("E3").Value = '  -1.41%  '

$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '631.13'
$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.75'
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("D7").Value = '3.768.70'
$ws.Range("E7").Value = '  -1.34%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("E10").Value = '  -1.98%  '

$ws.Range("E11").Value = '  +0.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.76'
$ws.Range("E12").Value = '  +1.45%  '

$ws.Range("E13").Value = '  -4.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.77'
$ws.Range("E14").Value = '  -3.67%  '

$ws.Range("D15").Value = '4.404.50'
$ws.Range("E15").Value = '  -1.46%  '

$ws.Range("D16").Value = '3.769.01'
$ws.Range("E16").Value = '  -1.08%  '

$ws.Range("D17").Value = '68.985.34'
$ws.Range("E17").Value = '  -0.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.64'
$ws.Range("E18").Value = '  -3.53%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.01'
$ws.Range("E20").Value = '  -1.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '462.50'
$ws.Range("E21").Value = '  -1.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.51'
$ws.Range("E22").Value = '  -2.16%  '

$ws.Range("E23").Value = '  -1.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.17'
$ws.Range("E24").Value = '  -2.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000144'
$ws.Range("E25").Value = '  -5.91%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.11'
$ws.Range("E26").Value = '  +0.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.12'
$ws.Range("E27").Value = '  -1.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.12'
$ws.Range("E28").Value = '  +0.23%  '

$ws.Range("E29").Value = '  +0.02%  '

$ws.Range("D30").Value = '3.925.20'
$ws.Range("E30").Value = '  -1.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.28'
$ws.Range("E31").Value = '  +0.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.67'
$ws.Range("E32").Value = '  -0.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.08'
$ws.Range("E33").Value = '  -3.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '0.176'
$ws.Range("E34").Value = '  +18.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '28.41'
$ws.Range("E35").Value = '  -2.80%  '

$ws.Range("E36").Value = '  +0.10%  '

$ws.Range("D37").Value = '3.722.50'
$ws.Range("E37").Value = '  -1.36%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.90'
$ws.Range("E38").Value = '  -2.65%  '

$ws.Range("E39").Value = '  -1.50%  '

$ws.Range("E40").Value = '  -1.76%  '

$ws.Range("E41").Value = '  -2.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = '0.963'
$ws.Range("E42").Value = '  -1.58%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.10%  '

$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = '1.98'
$ws.Range("E45").Value = '  +4.40%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").Value = '156.69'
$ws.Range("E46").Value = '  +0.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.43'
$ws.Range("E47").Value = '  +0.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").Value = '43.29'
$ws.Range("E48").Value = '  +0.65%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = '46.92'
$ws.Range("E49").Value = '  +0.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.294'
$ws.Range("E50").Value = '  -2.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.35'
$ws.Range("E51").Value = '  -1.67%  '
